$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the cells we are about to write retain their original plain-text
# representation (e.g. "10", "308.82", "-4.05%") instead of being
# auto-converted to numbers/percentages by Excel.
$ws.Range("D2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "308.82"
$ws.Range("E2").Value = "-4.05%"
$ws.Range("G2").Value = "10"

# Row 3
$ws.Range("D3").Value = "40.01"
$ws.Range("E3").Value = "-6.07%"
$ws.Range("G3").Value = "10"

# Row 4
$ws.Range("D4").Value = "5.117"
$ws.Range("E4").Value = "-1.44%"
$ws.Range("G4").Value = "10"

# Row 5
$ws.Range("D5").Value = "0.07735"
$ws.Range("E5").Value = "-5.60%"
$ws.Range("G5").Value = "10"

# Row 6
$ws.Range("D6").Value = "4.260"
$ws.Range("E6").Value = "-0.71%"
$ws.Range("G6").Value = "10"

# Row 7
$ws.Range("D7").Value = "1.596"
$ws.Range("E7").Value = "-11.29%"
$ws.Range("G7").Value = "10"

# Row 8
$ws.Range("D8").Value = "0.8823"
$ws.Range("E8").Value = "-5.29%"
$ws.Range("G8").Value = "10"

# Row 9
$ws.Range("D9").Value = "0.09847"
$ws.Range("E9").Value = "-11.29%"
$ws.Range("G9").Value = "10"

# Row 10
$ws.Range("D10").Value = "0.1749"
$ws.Range("E10").Value = "-6.49%"
$ws.Range("G10").Value = "10"

# Row 11
$ws.Range("D11").Value = "0.09036"
$ws.Range("E11").Value = "-5.02%"
$ws.Range("G11").Value = "10"

# Row 12
$ws.Range("D12").Value = "0.04440"
$ws.Range("E12").Value = "-5.11%"
$ws.Range("G12").Value = "10"

# Row 13
$ws.Range("E13").Value = "-0.28%"
$ws.Range("G13").Value = "10"

# Row 14
$ws.Range("D14").Value = "0.001263"
$ws.Range("E14").Value = "-3.09%"
$ws.Range("G14").Value = "10"

# Row 15
$ws.Range("D15").Value = "0.005824"
$ws.Range("E15").Value = "2.57%"
$ws.Range("G15").Value = "10"

# Row 16
$ws.Range("E16").Value = "2,412.72%"
$ws.Range("G16").Value = "10"

# Row 17
$ws.Range("D17").Value = "3.355"
$ws.Range("E17").Value = "-0.25%"
$ws.Range("G17").Value = "10"

# Row 18
$ws.Range("E18").Value = "-4.06%"
$ws.Range("G18").Value = "10"

# Row 19
$ws.Range("D19").Value = "0.3274"
$ws.Range("E19").Value = "-2.96%"
$ws.Range("G19").Value = "10"

# Row 20
$ws.Range("D20").Value = "7.065"
$ws.Range("E20").Value = "-4.80%"
$ws.Range("G20").Value = "10"

# Row 21
$ws.Range("D21").Value = "0.1342"
$ws.Range("E21").Value = "-3.28%"
$ws.Range("G21").Value = "10"

# Row 22
$ws.Range("D22").Value = "0.2790"
$ws.Range("E22").Value = "11.80%"
$ws.Range("G22").Value = "10"

# Row 23
$ws.Range("D23").Value = "0.04133"
$ws.Range("E23").Value = "-0.38%"
$ws.Range("G23").Value = "10"

# Row 24
$ws.Range("D24").Value = "0.001201"
$ws.Range("E24").Value = "-3.59%"
$ws.Range("G24").Value = "10"

# Row 25
$ws.Range("D25").Value = "0.004096"
$ws.Range("E25").Value = "-6.09%"
$ws.Range("G25").Value = "10"

# Row 26
$ws.Range("D26").Value = "0.0001302"
$ws.Range("E26").Value = "8.38%"
$ws.Range("G26").Value = "10"

# Row 27
$ws.Range("G27").Value = "10"

# Row 28
$ws.Range("G28").Value = "10"

# Row 29
$ws.Range("G29").Value = "10"

# Row 30
$ws.Range("G30").Value = "10"

# Row 31
$ws.Range("G31").Value = "10"

# Row 32
$ws.Range("G32").Value = "10"

# Row 33
$ws.Range("G33").Value = "10"

# Row 34
$ws.Range("G34").Value = "10"

# Row 35
$ws.Range("G35").Value = "10"

# Row 36
$ws.Range("G36").Value = "10"

# Row 37
$ws.Range("G37").Value = "10"

# Row 38
$ws.Range("D38").Value = "0.02357"
$ws.Range("E38").Value = "-14.63%"
$ws.Range("G38").Value = "10"

# Row 39
$ws.Range("D39").Value = "0.05218"
$ws.Range("E39").Value = "-6.96%"
$ws.Range("G39").Value = "10"

# Row 40
$ws.Range("D40").Value = "0.007923"
$ws.Range("E40").Value = "-1.00%"
$ws.Range("G40").Value = "10"

# Row 41
$ws.Range("D41").Value = "0.1326"
$ws.Range("E41").Value = "-5.20%"
$ws.Range("G41").Value = "10"

# Row 42
$ws.Range("D42").Value = "0.006508"
$ws.Range("E42").Value = "-0.65%"
$ws.Range("G42").Value = "10"

# Row 43
$ws.Range("D43").Value = "0.001952"
$ws.Range("E43").Value = "-6.45%"
$ws.Range("G43").Value = "10"

# Row 44
$ws.Range("D44").Value = "0.008753"
$ws.Range("E44").Value = "4.90%"
$ws.Range("G44").Value = "10"

# Row 45
$ws.Range("D45").Value = "0.3343"
$ws.Range("E45").Value = "-4.53%"
$ws.Range("G45").Value = "10"

# Row 46
$ws.Range("D46").Value = "0.00006561"
$ws.Range("E46").Value = "-5.72%"
$ws.Range("G46").Value = "10"

# Row 47
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").Value = "0.07%"
$ws.Range("G47").Value = "10"

# Row 48
$ws.Range("E48").Value = "98.31%"
$ws.Range("G48").Value = "10"

# Row 49
$ws.Range("D49").Value = "0.003472"
$ws.Range("E49").Value = "-0.18%"
$ws.Range("G49").Value = "10"

# Row 50
$ws.Range("D50").Value = "0.00002103"
$ws.Range("E50").Value = "0.07%"
$ws.Range("G50").Value = "10"

# Row 51
$ws.Range("D51").Value = "0.0002003"
$ws.Range("E51").Value = "0.07%"
$ws.Range("G51").Value = "10"
